# RegressionTestSuite.xlsx - "Checked in New Code"
# Adds 4 new regression test rows (TC_Exel_005 .. TC_Exel_008) following the
# same pattern as the existing Exelon/PHI/Web rows, normalizes the border
# styling of rows 5-6 to match the rest of the table, and leaves the
# selection where the author last left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize formatting of rows 5-6 -------------------------------------
# In the original file rows 5-6 used a slightly different cell style
# (left/right border only). Copy the format used by the rows above (full
# thin box border) onto rows 5-6 so the whole table is consistently styled.
$ws.Range("A4:G4").Copy() | Out-Null
$ws.Range("A5:G6").PasteSpecial(-4122) | Out-Null

# --- Add four new rows (7-10) with the same formatting ---------------------
$ws.Range("A4:G4").Copy() | Out-Null
$ws.Range("A7:G10").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Fill in the new test case data ----------------------------------------
$newRows = @(
    @(6, "Exelon", "PHI", "Web", "Firfox", "TC_Exel_005", "Y"),
    @(7, "Exelon", "PHI", "Web", "Firfox", "TC_Exel_006", "Y"),
    @(8, "Exelon", "PHI", "Web", "Firfox", "TC_Exel_007", "Y"),
    @(9, "Exelon", "PHI", "Web", "Firfox", "TC_Exel_008", "Y")
)

$r = 7
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

# --- Restore the author's final selection -----------------------------------
$ws.Range("E12").Select() | Out-Null
